# Correction of MDV status of SID
#
# The Maldives (MDV, row 18) was incorrectly flagged as not being a SIDS
# (Small Island Developing State). Fix the SIDS column (D) for MDV, then
# (re-)apply the AutoFilter on the CPC table (A1:H31) so only the rows
# flagged as SIDS (column D = 1) are shown, and leave the selection on
# the cell below the corrected one, matching the author's workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPC")

# Correct the SIDS flag for MDV (Maldives) from 0 to 1
$ws.Range("D18").Value = 1

# Re-apply the AutoFilter over the table, filtering column D (SIDS, the
# 4th column of the range) to only show the rows where SIDS = 1
$rng = $ws.Range("A1:H31")
$rng.AutoFilter(4, @("1"), 7)

# Leave the active selection on D19, as left by the editor
$ws.Range("D19").Select()
